$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are kept as text, matching the source
# data which stores prices like "1.001" or "27.468.20" as literal strings
# rather than numbers (avoids Excel auto-converting "1.001" -> 1.001 etc.)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.468.20'
$ws.Range('E2').Value = '  -3.16%  '
$ws.Range('D3').Value = '1.751.12'
$ws.Range('E3').Value = '  -3.69%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('D5').Value = '322.20'
$ws.Range('E5').Value = '  -2.30%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').Value = '0.4256'
$ws.Range('E7').Value = '  -4.39%  '
$ws.Range('E8').Value = '  -3.30%  '
$ws.Range('D9').Value = '42.31'
$ws.Range('E9').Value = '  -6.08%  '
$ws.Range('D10').Value = '0.07465'
$ws.Range('E10').Value = '  -3.71%  '
$ws.Range('D11').Value = '1.096'
$ws.Range('E11').Value = '  -3.66%  '
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').Value = '20.67'
$ws.Range('E13').Value = '  -7.04%  '
$ws.Range('D14').Value = '6.022'
$ws.Range('E14').Value = '  -4.94%  '
$ws.Range('D15').Value = '7.199'
$ws.Range('E15').Value = '  -5.70%  '
$ws.Range('D16').Value = '1.747.95'
$ws.Range('E16').Value = '  -5.30%  '
$ws.Range('D17').Value = '93.56'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range('D19').Value = '0.06346'
$ws.Range('E19').Value = '  -2.85%  '
$ws.Range('D20').Value = '0.9998'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '17.00'
$ws.Range('E21').Value = '  -3.35%  '
$ws.Range('D22').Value = '5.891'
$ws.Range('E22').Value = '  -5.83%  '
$ws.Range('D23').Value = '27.512.92'
$ws.Range('E23').Value = '  -3.13%  '
$ws.Range('D24').Value = '11.20'
$ws.Range('E24').Value = '  -4.37%  '
$ws.Range('D25').Value = '2.091'
$ws.Range('E25').Value = '  -4.07%  '
$ws.Range('D26').Value = '162.13'
$ws.Range('E26').Value = '  +3.32%  '
$ws.Range('D27').Value = '20.22'
$ws.Range('E27').Value = '  -3.04%  '
$ws.Range('D28').Value = '1.943.04'
$ws.Range('D29').Value = '2.135'
$ws.Range('E29').Value = '  -8.15%  '
$ws.Range('D30').Value = '123.93'
$ws.Range('E30').Value = '  -3.90%  '
$ws.Range('D31').Value = '1.099'
$ws.Range('E31').Value = '  -9.36%  '
$ws.Range('D32').Value = '3.665'
$ws.Range('E32').Value = '  +0.84%  '
$ws.Range('D33').Value = '5.558'
$ws.Range('E33').Value = '  -6.51%  '
$ws.Range('D34').Value = '0.08887'
$ws.Range('E34').Value = '  -4.03%  '
$ws.Range('D35').Value = '12.22'
$ws.Range('E35').Value = '  -8.08%  '
$ws.Range('D36').Value = '0.02283'
$ws.Range('E36').Value = '  -3.46%  '
$ws.Range('E37').Value = '  -4.90%  '
$ws.Range('D38').Value = '0.05996'
$ws.Range('E38').Value = '  -4.08%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = '4.955'
$ws.Range('E39').Value = '  -4.94%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.6306'
$ws.Range('E40').Value = '  -4.83%  '
$ws.Range('D41').Value = '1.187'
$ws.Range('E41').Value = '  -1.36%  '
$ws.Range('B42').Value = 'Frax'
$ws.Range('C42').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D42').Value = '0.9999'
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '7.919'
$ws.Range('E43').Value = '  -3.31%  '
$ws.Range('D44').Value = '1.390'
$ws.Range('E44').Value = '  -3.35%  '
$ws.Range('D45').Value = '13.42'
$ws.Range('E45').Value = '  -3.72%  '
$ws.Range('D46').Value = '0.5870'
$ws.Range('E46').Value = '  -4.72%  '
$ws.Range('E47').Value = '  -2.32%  '
$ws.Range('D48').Value = '123.35'
$ws.Range('E48').Value = '  -3.06%  '
$ws.Range('D49').Value = '1.966'
$ws.Range('E49').Value = '  -4.06%  '
$ws.Range('D50').Value = '1.164'
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('D51').Value = '0.06830'
$ws.Range('E51').Value = '  -2.62%  '
